$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab to reflect the new "through" date
$ws.Name = "Through 2022-06-28"

# Update the label for the June row (shared string)
$ws.Range("A7").Value = "June (through 06-28)"

# Update the June row values (row 7) - columns C through I
$ws.Range("C7").Value = 37
$ws.Range("D7").Value = 70
$ws.Range("E7").Value = 54
$ws.Range("F7").Value = 43
$ws.Range("G7").Value = 108
$ws.Range("H7").Value = 115
$ws.Range("I7").Value = 133

# Update the Total row values (row 8) - columns C through I
$ws.Range("C8").Value = 246
$ws.Range("D8").Value = 386
$ws.Range("E8").Value = 349
$ws.Range("F8").Value = 247
$ws.Range("G8").Value = 466
$ws.Range("H8").Value = 746
$ws.Range("I8").Value = 796
